$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything in the previously used range (A1:H8) first so that
# rows/columns beyond the new D4 boundary are removed entirely.
$ws.Range("A1:H8").Clear()

# Re-write the surviving cells with their new (shorter) "kodepython" spelling.
$ws.Range("A1").Value = "k"

$ws.Range("A2").Value = "o"
$ws.Range("B2").Value = "d"

$ws.Range("A3").Value = "e"
$ws.Range("B3").Value = "p"
$ws.Range("C3").Value = "y"

$ws.Range("A4").Value = "t"
$ws.Range("B4").Value = "h"
$ws.Range("C4").Value = "o"
$ws.Range("D4").Value = "n"
